# Update NATMI LR-pairs sheet (Efna3-Epha4) with newly recomputed TPM-based
# ligand/receptor expression and specificity values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Target cluster: ECs) ---
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.370913
$ws.Range("H2").Value = 1.112739
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 10.23061133333333
$ws.Range("N2").Value = 30.691834
$ws.Range("O2").Value = 0.4855635428718841
$ws.Range("P2").Value = 0.4855635428718841
$ws.Range("Q2").Value = 3.794666741480667
$ws.Range("R2").Value = 34.152000673326
$ws.Range("S2").Value = 0.4855635428718841
$ws.Range("T2").Value = 0.4855635428718841

# --- Row 3 (Target cluster: FAPs) ---
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.370913
$ws.Range("H3").Value = 1.112739
$ws.Range("O3").Value = 0.4164864079521221
$ws.Range("P3").Value = 0.4164864079521222
$ws.Range("Q3").Value = 3.254830688455666
$ws.Range("R3").Value = 29.293476196101
$ws.Range("S3").Value = 0.4164864079521221
$ws.Range("T3").Value = 0.4164864079521222

# --- Row 4 (Target cluster: MuSCs) ---
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.370913
$ws.Range("H4").Value = 1.112739
$ws.Range("M4").Value = 2.034752
$ws.Range("N4").Value = 6.104255999999999
$ws.Range("O4").Value = 0.09657305490303886
$ws.Range("P4").Value = 0.09657305490303887
$ws.Range("Q4").Value = 0.7547159685759999
$ws.Range("R4").Value = 6.792443717183999
$ws.Range("S4").Value = 0.09657305490303886
$ws.Range("T4").Value = 0.09657305490303887

# --- Row 5 (Target cluster: Resolving-Mac) ---
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.370913
$ws.Range("H5").Value = 1.112739
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.02901266666666667
$ws.Range("N5").Value = 0.087038
$ws.Range("O5").Value = 0.001376994272954919
$ws.Range("P5").Value = 0.001376994272954919
$ws.Range("Q5").Value = 0.01076117523133333
$ws.Range("R5").Value = 0.096850577082
$ws.Range("S5").Value = 0.001376994272954919
$ws.Range("T5").Value = 0.001376994272954919
